$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "أحمد سعيد"
$ws.Range("A2").Value = "ahmed 2"
$ws.Range("A3").Value = "ahmed 3"
$ws.Range("A4").Value = "ahmed 4"

$ws.Range("A2").Select()
